$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# New header cells (date / legislator_name / legislator_id)
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# Pre-format column H as text so the date-like string "2013-12-19" is kept
# literal instead of being auto-converted into a date serial number.
$ws.Range("H2:H8").NumberFormat = "@"

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = "2013-12-19"
    $ws.Cells.Item($r, 9).Value = "楊玉欣"
    $ws.Cells.Item($r, 10).Value = 1757
}
